# Remove trailing whitespace from the "MCA temporarily removed" caveat
# messages in column B (dataText sheet). The author's commit message
# states: "removed whitespace which should remove double full stop on
# MCA message." Only the trailing space before the closing of each
# string is removed; the rest of the text (including the internal
# double space after the date prefix) is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = "Jul 2023-Jun 2024.  North East and North of Tyne MCA have temporarily been removed, on account of North of Tyne being integrated into North East from now on."
$ws.Cells.Item(3, 2).Value = "Jul 2023-Jun 2024.  North East and North of Tyne MCA have temporarily been removed, on account of North of Tyne being integrated into North East from now on."
$ws.Cells.Item(4, 2).Value = "Jul 2023-Jun 2024.  North East and North of Tyne MCA have temporarily been removed, on account of North of Tyne being integrated into North East from now on."
$ws.Cells.Item(5, 2).Value = "Jul 2023-Jun 2024.  North East and North of Tyne MCA have temporarily been removed, on account of North of Tyne being integrated into North East from now on."
$ws.Cells.Item(6, 2).Value = "Jul 2023-Jun 2024.  North East and North of Tyne MCA have temporarily been removed, on account of North of Tyne being integrated into North East from now on."
$ws.Cells.Item(7, 2).Value = "Jul 2023-Jun 2024.  North East and North of Tyne MCA have temporarily been removed, on account of North of Tyne being integrated into North East from now on."
$ws.Cells.Item(8, 2).Value = "Jul 2023-Jun 2024.  North East and North of Tyne MCA have temporarily been removed, on account of North of Tyne being integrated into North East from now on."
$ws.Cells.Item(9, 2).Value = "Jul 2023-Jun 2024.  North East and North of Tyne MCA have temporarily been removed, on account of North of Tyne being integrated into North East from now on."
$ws.Cells.Item(10, 2).Value = "Nov 2024 data. North East and North of Tyne MCA have temporarily been removed, on account of North of Tyne being integrated into North East from now on."
$ws.Cells.Item(11, 2).Value = "Mar 2024 data.  North East and North of Tyne MCA have temporarily been removed, on account of North of Tyne being integrated into North East from now on."
$ws.Cells.Item(12, 2).Value = "Dec 2021 - Dec 2022 data.  North East and North of Tyne MCA have temporarily been removed, on account of North of Tyne being integrated into North East from now on."
$ws.Cells.Item(13, 2).Value = "Dec 2021 - Dec 2022 data.  North East and North of Tyne MCA have temporarily been removed, on account of North of Tyne being integrated into North East from now on."
$ws.Cells.Item(14, 2).Value = "AY22/23 data.  North East and North of Tyne MCA have temporarily been removed, on account of North of Tyne being integrated into North East from now on."
$ws.Cells.Item(15, 2).Value = "AY22/23 data.  North East and North of Tyne MCA have temporarily been removed, on account of North of Tyne being integrated into North East from now on."
$ws.Cells.Item(16, 2).Value = "AY22/23 data.  North East and North of Tyne MCA have temporarily been removed, on account of North of Tyne being integrated into North East from now on."
$ws.Cells.Item(17, 2).Value = "AY22/23 data.  North East and North of Tyne MCA have temporarily been removed, on account of North of Tyne being integrated into North East from now on."
$ws.Cells.Item(18, 2).Value = "Jan-Dec 2023 data.  North East and North of Tyne MCA have temporarily been removed, on account of North of Tyne being integrated into North East from now on."
$ws.Cells.Item(19, 2).Value = "Jan-Dec 2023 data.  North East and North of Tyne MCA have temporarily been removed, on account of North of Tyne being integrated into North East from now on."
$ws.Cells.Item(20, 2).Value = "AY22/23 data.  North East and North of Tyne MCA have temporarily been removed, on account of North of Tyne being integrated into North East from now on."
$ws.Cells.Item(21, 2).Value = "AY22/23 data.  North East and North of Tyne MCA have temporarily been removed, on account of North of Tyne being integrated into North East from now on."
$ws.Cells.Item(22, 2).Value = "Growth from 2023 to 2035.  North East and North of Tyne MCA have temporarily been removed, on account of North of Tyne being integrated into North East from now on."
$ws.Cells.Item(23, 2).Value = "AY22/23 data.  North East and North of Tyne MCA have temporarily been removed, on account of North of Tyne being integrated into North East from now on."
$ws.Cells.Item(24, 2).Value = "AY22/23 data.  North East and North of Tyne MCA have temporarily been removed, on account of North of Tyne being integrated into North East from now on."
